$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New schedule data for rows 2..11 (letni semestr 2024 info kombi) ---
# Columns: A Den | B Od | C Do | D Predmet | E Akce | F Mistnost | G Vyucujici | H Omezeni | I Kapacita

$data = @(
  @("So 17.02.2024","08:45","16:45","EIE87E Systémová integrace - INFON4 komb.","Přednáška","AII","Tyrychtr Jan","2n-infonp"),
  @("Ne 18.02.2024","08:45","16:45","EIEB7E Gramatiky a jazyky - INFON4 komb.","Přednáška","PEF C11","Merunka Vojtěch","2n-infonp"),
  @("Pá 08.03.2024","12:15","20:15","ETEW8E Zpracování velkých dat - INFON4v Kombi.","Přednáška","TI","Masner Jan","2n-infonp"),
  @("So 09.03.2024","08:45","16:45","ETE90E Internetové technologie - server side - INFON4k komb.","Přednáška","PEF C11","Lohr Václav","2n-infonp"),
  @("Ne 10.03.2024","08:45","16:45","EJE98E Právní aspekty IT INFON4 DS","Přednáška","PEF C11","Reichert Michal","2n-infonp"),
  @("Pá 05.04.2024","13:00","21:00","EIEB7E Gramatiky a jazyky - INFON4 komb.","Cvičení","PEF C11","Merunka Vojtěch","2n-infonp"),
  @("So 06.04.2024","08:45","16:45","ETE90E Internetové technologie - server side - INFON4k komb. (1)","Cvičení","PEF C11","Lohr Václav","2n-infonp"),
  @("Ne 07.04.2024","08:45","16:45","EIE87E Systémová integrace - INFON4 komb.","Cvičení","PEF EII","Tyrychtr Jan","2n-infonp"),
  @("Pá 26.04.2024","13:05","21:05","ETEW8E Zpracování velkých dat - INFON4v Kombi.","Cvičení","AIII","Masner Jan","2n-infonp"),
  @("Ne 28.04.2024","08:45","16:45","EJE98E Právní aspekty IT INFON4 DS","Cvičení","PEF EII","Reichert Michal","2n-infonp")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $ws.Cells.Item($r, 9).ClearContents()
  $r = $r + 1
}

# Rows 12 and 13 from the old schedule no longer exist - remove them
$ws.Rows("12:13").Delete()

# Column width tweaks (D wider, H narrower) to fit the new (longer) subject names
$ws.Columns("D").ColumnWidth = 66.83333333333333
$ws.Columns("H").ColumnWidth = 11.833333333333334
